# WorkReport.xlsx update: lokalizace, export do Google Docs
# Adds three new work-log rows (62-64) to Sheet1, extends the SUM formula
# range, and updates the selection to reflect the next empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Copy date formatting (style) from the last existing date cell (D61)
# into the three new date cells, then overwrite their values. ---
$ws.Range("D61").Copy($ws.Range("D62"))
$ws.Range("D61").Copy($ws.Range("D63"))
$ws.Range("D61").Copy($ws.Range("D64"))

# --- Row 62 ---
$ws.Range("B62").Value = "Michani Choice items, Lokalizace"
$ws.Range("C62").Value = 2
$ws.Range("D62").Value = 40992

# --- Row 63 ---
$ws.Range("B63").Value = "Lokalizace, styly, Export do Google Docs"
$ws.Range("C63").Value = 5
$ws.Range("D63").Value = 40993

# --- Row 64 ---
$ws.Range("B64").Value = "Export do CSV, Google Docs"
$ws.Range("C64").Value = 2
$ws.Range("D64").Value = 40995

# --- Extend the total-hours SUM formula range in C3 to cover new rows ---
$ws.Range("C3").Formula = "=SUM(C4:C567)"

# --- Update the selected / active cell to the next empty row (B65) ---
[void]$ws.Range("B65").Select()
